$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-08-12"

# Update the header label in I1
$ws.Range("I1").Value = "2022 (through 08-12)"

# Update the August row (row 9) value for 2022 column (I)
$ws.Range("I9").Value = 66

# Update the Total row (row 14) value for 2022 column (I)
$ws.Range("I14").Value = 1036
